# Adds the new "UDT/Name" and "UDT/Identifier(...)" terms to the Terms
# sheet and fills in the English ("en") labels/descriptions for them.
# The other language sheets (es, fr, pt, si) pick the new terms up
# automatically because their A7:A15 cells already hold "=Terms!A7" style
# formulas that were previously resolving to 0 (Terms!A7:A15 was blank).

$wb = $excel.ActiveWorkbook

$terms = $wb.Worksheets.Item("Terms")
$en    = $wb.Worksheets.Item("en")

# ---------------------------------------------------------------------
# 1. Terms sheet - new rows 7-15 (term path / ccts type / datatype / optional)
# ---------------------------------------------------------------------
$termsRows = @(
    @{ Row=7;  A="UDT/Name";                        B="ccts-cct:TextType";       C="xsd:string";             D="00" },
    @{ Row=8;  A="UDT/Identifier";                   B="ccts-cct:IdentifierType"; C="xsd:normalizedString";   D="00" },
    @{ Row=9;  A="UDT/Identifier/@schemeID";         B="ccts-cct:IdentifierType"; C="xsd:normalizedString";   D="00" },
    @{ Row=10; A="UDT/Identifier/@schemeName";       B="ccts-cct:IdentifierType"; C="xsd:string";             D="00" },
    @{ Row=11; A="UDT/Identifier/@schemeAgencyID";   B="ccts-cct:IdentifierType"; C="xsd:normalizedString";   D="00" },
    @{ Row=12; A="UDT/Identifier/@schemeAgencyName"; B="ccts-cct:IdentifierType"; C="xsd:string";             D="00" },
    @{ Row=13; A="UDT/Identifier/@schemeVersionID";  B="ccts-cct:IdentifierType"; C="xsd:normalizedString";   D="00" },
    @{ Row=14; A="UDT/Identifier/@schemeDataURI";    B="ccts-cct:IdentifierType"; C="xsd:anyURI";             D="00" },
    @{ Row=15; A="UDT/Identifier/@schemeURI";        B="ccts-cct:IdentifierType"; C="xsd:anyURI";             D="00" }
)

# Column D ("Optional") on this sheet is stored/displayed as text (e.g. "00"),
# mirror the existing D2:D6 cells by forcing a text number format first so
# Excel doesn't coerce "00" into the number 0.
$terms.Range("D7:D15").NumberFormat = "@"

foreach ($r in $termsRows) {
    $terms.Cells.Item($r.Row, 1).Value = $r.A
    $terms.Cells.Item($r.Row, 2).Value = $r.B
    $terms.Cells.Item($r.Row, 3).Value = $r.C
    $terms.Cells.Item($r.Row, 4).Value = $r.D
}

# The existing list validation on D2:D6 needs to cover the new rows too -
# recreate it across D2:D15 so it collapses back into a single sqref.
$terms.Range("D2:D15").Validation.Delete()
$terms.Range("D2:D15").Validation.Add(3, 1, 1, '"00,01,10,11"')

# ---------------------------------------------------------------------
# 2. "en" sheet - English label/description for the same new rows.
#    Column A is already driven by "=Terms!Ann" formulas and will pick
#    up the new values automatically once Terms! has data.
# ---------------------------------------------------------------------
$enRows = @(
    @{ Row=7;  B="Name";               C="A character string that constitutes the distinctive designation of a person, place, thing or concept." },
    @{ Row=8;  B="Identifier";         C="A character string to identify and distinguish uniquely, one instance of an object in an identification scheme from all other objects in the same scheme together with relevant supplementary information." },
    @{ Row=9;  B="Scheme ID";          C="The identification of the identification scheme." },
    @{ Row=10; B="Scheme Name";        C="The name of the identification scheme." },
    @{ Row=11; B="Scheme Agency ID";   C="The identification of the agency that maintains the identification scheme." },
    @{ Row=12; B="Scheme Agency Name"; C="The name of the agency that maintains the identification scheme." },
    @{ Row=13; B="Scheme Version ID";  C="The version of the identification scheme." },
    @{ Row=14; B="Scheme Data URI";    C="The Uniform Resource Identifier that identifies where the identification scheme data is located." },
    @{ Row=15; B="Scheme URI";         C="The Uniform Resource Identifier that identifies where the identification scheme is located." }
)

foreach ($r in $enRows) {
    $en.Cells.Item($r.Row, 2).Value = $r.B
    $en.Cells.Item($r.Row, 3).Value = $r.C
    $en.Cells.Item($r.Row, 5).Value = $true
}

# Row 8 wraps onto two lines in the authored workbook.
$en.Rows.Item(8).RowHeight = 28.8

# Column B on both sheets widens to fit the newly entered labels.
$terms.Columns.Item(2).ColumnWidth = 19.44140625
$en.Columns.Item(2).ColumnWidth = 24.33203125

# ---------------------------------------------------------------------
# 3. Selection bookkeeping: the author ended up with C8 selected on the
#    "en" sheet, while the Terms sheet stays the active tab.
# ---------------------------------------------------------------------
$en.Range("C8").Select()
$terms.Activate()
